# Auto-generated PowerShell Excel COM-interop script
# Applies cryptocurrency price/volume table updates per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D that would otherwise be auto-parsed as numbers by Excel need to be
# forced to Text format first, so the stored value stays an exact string (matching
# the original inline-string representation, e.g. "570.50" not 570.5).
$textCells = @(
    "D4", "D5", "D6", "D7", "D10", "D12", "D13", "D14", "D16", "D20",
    "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31",
    "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D43", "D44",
    "D46", "D47", "D49", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range("D2").Value = '66.938.30'
$ws.Range("E2").Value = '  +3.02%  '
$ws.Range("D3").Value = '3.440.26'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '570.50'
$ws.Range("E5").Value = '  +1.94%  '
$ws.Range("D6").Value = '184.40'
$ws.Range("E6").Value = '  +5.64%  '
$ws.Range("D7").Value = '0.634'
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("D8").Value = '3.432.76'
$ws.Range("E8").Value = '  +1.52%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = '0.177'
$ws.Range("E10").Value = '  +6.53%  '
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").Value = '55.33'
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("D13").Value = '0.0000281'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = '9.37'
$ws.Range("E14").Value = '  +2.87%  '
$ws.Range("D15").Value = '3.982.18'
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").Value = '18.53'
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").Value = '3.429.75'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '66.722.71'
$ws.Range("E19").Value = '  +2.84%  '
$ws.Range("D20").Value = '12.02'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("E21").Value = '  +1.76%  '
$ws.Range("D22").Value = '469.88'
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").Value = '4.98'
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").Value = '14.90'
$ws.Range("E24").Value = '  +9.67%  '
$ws.Range("D25").Value = '4.19'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").Value = '89.58'
$ws.Range("E26").Value = '  +3.37%  '
$ws.Range("D27").Value = '2.96'
$ws.Range("E27").Value = '  +0.53%  '
$ws.Range("D28").Value = '10.96'
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").Value = '8.92'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").Value = '31.49'
$ws.Range("E30").Value = '  +2.29%  '
$ws.Range("D31").Value = '6.98'
$ws.Range("E31").Value = '  +2.89%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").Value = '11.62'
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '62.93'
$ws.Range("E33").Value = '  +2.59%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = '581.34'
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("E35").Value = '  +1.86%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.148'
$ws.Range("E36").Value = '  +6.26%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '3.65'
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("D39").Value = '0.391'
$ws.Range("E39").Value = '  +5.94%  '
$ws.Range("D40").Value = '36.60'
$ws.Range("E40").Value = '  +2.98%  '
$ws.Range("D41").Value = '0.0₃0769'
$ws.Range("E41").Value = '  +3.71%  '
$ws.Range("D42").Value = '3.124.45'
$ws.Range("E42").Value = '  +1.28%  '
$ws.Range("D43").Value = '2.92'
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").Value = '2.62'
$ws.Range("E44").Value = '  +6.71%  '
$ws.Range("E45").Value = '  +2.65%  '
$ws.Range("D46").Value = '2.79'
$ws.Range("E46").Value = '  +20.24%  '
$ws.Range("D47").Value = '3.25'
$ws.Range("E47").Value = '  +3.18%  '
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").Value = '0.997'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("E50").Value = '  +2.65%  '
$ws.Range("D51").Value = '8.66'
$ws.Range("E51").Value = '  +4.64%  '
